{"js": "// Update the date line and every \"dividend\u00f7divisor=quotient, remainder\"\n// answer cell in the worksheet table. Replacements are applied\n// positionally (in document order) because some new values equal other\n// cells' old values, so a global text search/replace would misfire.\n\n// Build the ordered list of old->new text pairs exactly as they occur in\n// the document (title paragraph first, then each table cell in reading\n// order: row-by-row, left-to-right). Rows with no content are skipped\n// automatically because their paragraph text is empty.\nconst pairs = [\n  [\"2025-07-31 Thursday\", \"2025-08-01 Friday\"],\n  [\"640\u00f77=91, 3\", \"679\u00f75=135, 4\"],\n  [\"297\u00f78=37, 1\", \"562\u00f79=62, 4\"],\n  [\"399\u00f79=44, 3\", \"798\u00f76=133, 0\"],\n  [\"290\u00f75=58, 0\", \"388\u00f78=48, 4\"],\n  [\"732\u00f79=81, 3\", \"973\u00f74=243, 1\"],\n  [\"542\u00f76=90, 2\", \"702\u00f77=100, 2\"],\n  [\"254\u00f73=84, 2\", \"328\u00f74=82, 0\"],\n  [\"912\u00f75=182, 2\", \"195\u00f79=21, 6\"],\n  [\"103\u00f78=12, 7\", \"625\u00f75=125, 0\"],\n  [\"132\u00f73=44, 0\", \"650\u00f77=92, 6\"],\n  [\"654\u00f79=72, 6\", \"657\u00f79=73, 0\"],\n  [\"755\u00f72=377, 1\", \"674\u00f73=224, 2\"],\n  [\"638\u00f72=319, 0\", \"921\u00f76=153, 3\"],\n  [\"594\u00f77=84, 6\", \"640\u00f77=91, 3\"],\n  [\"885\u00f77=126, 3\", \"351\u00f73=117, 0\"],\n  [\"171\u00f72=85, 1\", \"903\u00f72=451, 1\"],\n  [\"335\u00f73=111, 2\", \"697\u00f79=77, 4\"],\n  [\"814\u00f79=90, 4\", \"396\u00f72=198, 0\"],\n  [\"769\u00f72=384, 1\", \"932\u00f72=466, 0\"],\n  [\"512\u00f77=73, 1\", \"285\u00f75=57, 0\"],\n  [\"537\u00f73=179, 0\", \"178\u00f78=22, 2\"],\n  [\"216\u00f79=24, 0\", \"372\u00f78=46, 4\"],\n  [\"634\u00f77=90, 4\", \"905\u00f75=181, 0\"],\n  [\"700\u00f78=87, 4\", \"187\u00f77=26, 5\"],\n  [\"737\u00f77=105, 2\", \"119\u00f76=19, 5\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"\") {\n    continue; // blank spacer row/paragraph, nothing to replace\n  }\n  if (idx >= pairs.length) {\n    break;\n  }\n  const [oldText, newText] = pairs[idx];\n  if (text !== oldText) {\n    throw new Error(\n      \"Unexpected paragraph text at position \" +\n        idx +\n        \": expected \" +\n        JSON.stringify(oldText) +\n        \" got \" +\n        JSON.stringify(text)\n    );\n  }\n  // Replace only the text content; formatting (run properties) on the\n  // paragraph's range is preserved by insertText(\"Replace\").\n  para.getRange().insertText(newText, Word.InsertLocation.replace);\n  idx++;\n}\n\nawait context.sync();\n\nif (idx !== pairs.length) {\n  throw new Error(\"Only replaced \" + idx + \" of \" + pairs.length + \" entries\");\n}\n", "ps1": "# Update the date line and every \"dividend\u00f7divisor=quotient, remainder\"\n# answer cell in the worksheet table. Replacements are applied\n# positionally (in document order) because some new values equal other\n# cells' old values, so a global Find/Replace would misfire (e.g. the\n# old text of cell 2 becomes the new text of cell 15).\n\n$d = $word.ActiveDocument\n\n# Ordered old->new text pairs exactly as they occur in the document:\n# the title paragraph first, then each non-empty table cell in reading\n# order (row-by-row, left-to-right). Empty cells are skipped.\n$pairs = @(\n    ,(\"2025-07-31 Thursday\", \"2025-08-01 Friday\")\n    ,(\"640\u00f77=91, 3\", \"679\u00f75=135, 4\")\n    ,(\"297\u00f78=37, 1\", \"562\u00f79=62, 4\")\n    ,(\"399\u00f79=44, 3\", \"798\u00f76=133, 0\")\n    ,(\"290\u00f75=58, 0\", \"388\u00f78=48, 4\")\n    ,(\"732\u00f79=81, 3\", \"973\u00f74=243, 1\")\n    ,(\"542\u00f76=90, 2\", \"702\u00f77=100, 2\")\n    ,(\"254\u00f73=84, 2\", \"328\u00f74=82, 0\")\n    ,(\"912\u00f75=182, 2\", \"195\u00f79=21, 6\")\n    ,(\"103\u00f78=12, 7\", \"625\u00f75=125, 0\")\n    ,(\"132\u00f73=44, 0\", \"650\u00f77=92, 6\")\n    ,(\"654\u00f79=72, 6\", \"657\u00f79=73, 0\")\n    ,(\"755\u00f72=377, 1\", \"674\u00f73=224, 2\")\n    ,(\"638\u00f72=319, 0\", \"921\u00f76=153, 3\")\n    ,(\"594\u00f77=84, 6\", \"640\u00f77=91, 3\")\n    ,(\"885\u00f77=126, 3\", \"351\u00f73=117, 0\")\n    ,(\"171\u00f72=85, 1\", \"903\u00f72=451, 1\")\n    ,(\"335\u00f73=111, 2\", \"697\u00f79=77, 4\")\n    ,(\"814\u00f79=90, 4\", \"396\u00f72=198, 0\")\n    ,(\"769\u00f72=384, 1\", \"932\u00f72=466, 0\")\n    ,(\"512\u00f77=73, 1\", \"285\u00f75=57, 0\")\n    ,(\"537\u00f73=179, 0\", \"178\u00f78=22, 2\")\n    ,(\"216\u00f79=24, 0\", \"372\u00f78=46, 4\")\n    ,(\"634\u00f77=90, 4\", \"905\u00f75=181, 0\")\n    ,(\"700\u00f78=87, 4\", \"187\u00f77=26, 5\")\n    ,(\"737\u00f77=105, 2\", \"119\u00f76=19, 5\")\n)\n\n$idx = 0\n\n# 1) Title paragraph (outside the table).\n$titleRange = $d.Paragraphs.Item(1).Range\n$titleText = $titleRange.Text.TrimEnd([char]13)\n$pair = $pairs[$idx]\n$expectedOld = $pair[0]\n$expectedNew = $pair[1]\nif ($titleText -ne $expectedOld) {\n    throw \"Unexpected title text: expected [$expectedOld] got [$titleText]\"\n}\n$titleRange.Text = $expectedNew\n$idx++\n\n# 2) Table cells, row-by-row / left-to-right; blank cells are skipped.\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($cellText.Length -eq 0) {\n            continue\n        }\n        if ($idx -ge $pairs.Count) {\n            throw \"Ran out of replacement pairs at row=$r col=$c\"\n        }\n        $pair = $pairs[$idx]\n        $expectedOld = $pair[0]\n        $expectedNew = $pair[1]\n        if ($cellText -ne $expectedOld) {\n            throw \"Unexpected cell text at row=$r col=$c : expected [$expectedOld] got [$cellText]\"\n        }\n        $t.Cell($r, $c).Range.Text = $expectedNew\n        $idx++\n    }\n}\n\nif ($idx -ne $pairs.Count) {\n    throw \"Only replaced $idx of $($pairs.Count) entries\"\n}\n\n\"done: replaced $idx entries\"\n"}
